$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in totals for the two previously-incomplete rows (BC817 / S13-24V-2C)
$ws.Range("I16").Formula = "=G16*H16"
$ws.Range("I17").Formula = "=G17*H17"

# New BOM line: USB ESD filter
$ws.Range("C18").Value = "USBUF02W6"
$ws.Range("D18").Value = "SOT323-6L"
$ws.Range("F18").Value = "Filtr: cyfrowy; line terminator; EMI,dolnoprzepustowy; SOT323-6L"
$ws.Range("B18").Value = "STM"
$ws.Range("E18").Value = "SMT"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1.0642
$ws.Range("I18").Formula = "=G18*H18"
$ws.Range("J18").Value = "https://www.tme.eu/pl/details/usbuf02w6/filtry-uklady-scalone/stmicroelectronics/"

$ws.Range("I23").Select() | Out-Null
